# Trade #8 closed at 2026-02-17 13:34:14 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly-closed trade #8 as a new
# row on both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.9    # Current Capital
$summary.Range("B4").Value = -0.1      # Total P&L $
$summary.Range("B5").Value = -0.25     # Total P&L %
$summary.Range("B6").Value = 8         # Total Trades
$summary.Range("B7").Value = 2         # Winning Trades
$summary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.90000000000001  # Capital
$status.Range("D4").Value = 8                  # Trades
$status.Range("E4").Value = -0.1               # P&L $
$status.Range("F4").Value = -0.1               # P&L %
$status.Range("G4").Value = 25                 # Win Rate %

# ---------------------------------------------------------------------
# Append the new closed trade (#8) to the trade-log sheets.
# Text-like values that could otherwise be misread as dates/times are
# written with a leading apostrophe so Excel keeps them as plain text
# (matching the other rows already on the sheet), then ClearFormats()
# strips the quote-prefix styling Excel applies for the apostrophe so
# the cell is left on the sheet's default (unstyled) format.
# ---------------------------------------------------------------------
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).ClearFormats()
}

function Add-TradeRow($ws) {
    $ws.Range("A9").Value = 8
    Set-TextCell $ws "B9" "2026-02-17"
    Set-TextCell $ws "C9" "13:34:08"
    $ws.Range("D9").Value = "MarketMaking"
    $ws.Range("E9").Value = "UP"
    $ws.Range("F9").Value = 0.03
    $ws.Range("G9").Value = 0.071577
    $ws.Range("H9").Value = "CLOSED"
    $ws.Range("I9").Value = 138.5895
    $ws.Range("J9").Value = 0.04
    $ws.Range("K9").Value = 99.90000000000001
    $ws.Range("L9").Value = 0
    $ws.Range("M9").Value = 0
    $ws.Range("N9").Value = 0.6
    $ws.Range("O9").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P9").Value = "early_exit"
    $ws.Range("Q9").Value = 0.1
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
